$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Invisible Digital Front: Can Cyber Attacks Shape Battlefield Events?"
$ws.Range("C2").Value = "Nadiya Kostyuk, Yuri M. Zhukov"
$ws.Range("D2").Value = "2019"
$ws.Range("E2").Value = "10.1177/0022002717737138"
$ws.Range("F2").Value = "Restricted"

# Row 3
$ws.Range("B3").Value = "Cyber scares and prophylactic policies: Crossnational evidence on the effect of cyberattacks on public support for surveillance"
$ws.Range("C3").Value = "Amelia C Arsenault, Sarah E Kreps, Keren LG Snider, Daphna Canetti"
$ws.Range("D3").Value = "2024"
$ws.Range("E3").Value = "10.1177/00223433241233960"
$ws.Range("G3").Value = 1

# Row 4
$ws.Range("B4").Value = "How the process of discovering cyberattacks biases our understanding of cybersecurity"
$ws.Range("C4").Value = "Harry Oppenheimer"
$ws.Range("E4").Value = "10.1177/00223433231217687"
$ws.Range("F4").Value = "Open Access"
$ws.Range("G4").Value = 1

# Row 5
$ws.Range("B5").Value = "On domains: Cyber and the practice of warfare"
$ws.Range("C5").Value = "Chris McGuffin, Paul Mitchell"
$ws.Range("D5").Value = "2014"
$ws.Range("E5").Value = "10.1177/0020702014540618"
$ws.Range("G5").Value = 1

# Row 6
$ws.Range("B6").Value = "Using network digital twins to improve cyber resilience of missions"
$ws.Range("C6").Value = "Rajive Bagrodia"
$ws.Range("D6").Value = "2023"
$ws.Range("E6").Value = "10.1177/15485129221131226"

# Row 7
$ws.Range("B7").Value = "Towards a Chronology of Robotic Art"
$ws.Range("C7").Value = "Eduardo Kac"
$ws.Range("D7").Value = "2001"
$ws.Range("E7").Value = "10.1177/135485650100700109"
$ws.Range("F7").Value = "Restricted"

# Row 8
$ws.Range("B8").Value = "Digital Assays Part II: Digital Protein and Cell Assays"
$ws.Range("C8").Value = "Amar S. Basu"
$ws.Range("D8").Value = "2017"
$ws.Range("E8").Value = "10.1177/2472630317705681"

# Row 9
$ws.Range("B9").Value = "Simplification and Linearization of Manipulator Dynamics by the Design of Inertia Distribution"
$ws.Range("C9").Value = "D.C.H. Yang, S.W. Tzeng"
$ws.Range("D9").Value = "1986"
$ws.Range("E9").Value = "10.1177/027836498600500307"
$ws.Range("G9").Value = 1

# Row 10
$ws.Range("B10").Value = "A novel ensemble learning approach for fault detection of sensor data in cyber-physical system"
$ws.Range("C10").Value = "Ramesh Sneka Nandhini, Ramanathan Lakshmanan"
$ws.Range("D10").Value = "2023"
$ws.Range("E10").Value = "10.3233/JIFS-235809"

# Row 11
$ws.Range("B11").Value = "Responding to Uncertainty: The Importance of Covertness in Support for Retaliation to Cyber and Kinetic Attacks"
$ws.Range("C11").Value = "Kathryn Hedgecock, Lauren Sukin"
$ws.Range("E11").Value = "10.1177/00220027231153580"
$ws.Range("F11").Value = "Open Access"
$ws.Range("G11").Value = 1
